$wb = $excel.ActiveWorkbook

# --- Sheet 1: "维修任务模板" (repair task template) ---
$ws1 = $wb.Worksheets.Item(1)

# Row 1 headers: replace the placeholder "?" runs with the real Chinese text.
$ws1.Range("A1").Value = "任务名称"
$ws1.Range("B1").Value = "任务类别"
$ws1.Range("C1").Value = "给分方式"
$ws1.Range("D1").Value = "单位积分"
$ws1.Range("E1").Value = "数量"
$ws1.Range("F1").Value = "积分规则"
$ws1.Range("G1").Value = "数量是否可修改"
$ws1.Range("H1").Value = "积分是否可修改"

# Row 2 sample data: a new "数量" (quantity) column is inserted at E2, which
# pushes the old E2/F2 contents one column to the right (F2/G2), and a new
# "积分是否可修改" value is added at H2.
$ws1.Range("E2").Value = 1
$ws1.Range("F2").Value = "备注说明"
$ws1.Range("G2").Value = "否"
$ws1.Range("H2").Value = "否"

# --- Sheet 2: "填写说明" (fill-in instructions) ---
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("A1").Value = "表头"
$ws2.Range("B1").Value = "填写说明"

$ws2.Range("A2").Value = "任务名称"
$ws2.Range("B2").Value = "必填，任务名称，不能重复。"

$ws2.Range("A3").Value = "任务类别"
$ws2.Range("B3").Value = "选填，自由文本，用于筛选区分。"

$ws2.Range("A4").Value = "给分方式"
$ws2.Range("B4").Value = "选填，填写 奖扣结合式 / 扣分项 / 奖分项。"

$ws2.Range("A5").Value = "单位积分"
$ws2.Range("B5").Value = "必填，数字，可为正/负/0，示例 10 / -5。"

$ws2.Range("A6").Value = "数量"
$ws2.Range("B6").Value = "选填，1-1000 的整数，默认 1。"

$ws2.Range("A7").Value = "积分规则"
$ws2.Range("B7").Value = "选填，备注说明。"

$ws2.Range("A8").Value = "数量是否可修改"
$ws2.Range("B8").Value = "选填，填写 是/否，默认 否。"

$ws2.Range("A9").Value = "积分是否可修改"
$ws2.Range("B9").Value = "选填，填写 是/否，默认 否。"
